$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format so numeric-looking strings
# (e.g. "309.98") are preserved as text and not converted to numbers.
$textCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "E27", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.095.60"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "2.301.44"
$ws.Range("E3").Value = "  -2.96%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "309.98"
$ws.Range("E5").Value = "  -6.24%  "
$ws.Range("D6").Value = "104.18"
$ws.Range("E6").Value = "  +4.80%  "
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  -4.14%  "
$ws.Range("D10").Value = "39.65"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").Value = "0.0908"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "8.24"
$ws.Range("E12").Value = "  -2.62%  "
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "0.963"
$ws.Range("E14").Value = "  -4.39%  "
$ws.Range("D15").Value = "15.36"
$ws.Range("E15").Value = "  -6.03%  "
$ws.Range("D16").Value = "2.651.71"
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").Value = "2.304.32"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").Value = "42.075.12"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  -5.22%  "
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").Value = "74.49"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").Value = "3.46"
$ws.Range("E22").Value = "  -8.31%  "
$ws.Range("D23").Value = "260.13"
$ws.Range("E23").Value = "  -4.53%  "
$ws.Range("D24").Value = "2.25"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").Value = "9.17"
$ws.Range("E25").Value = "  -5.60%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "10.91"
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("D29").Value = "22.83"
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "35.56"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "164.14"
$ws.Range("E31").Value = "  -6.56%  "
$ws.Range("D32").Value = "0.0882"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("E33").Value = "  -6.53%  "
$ws.Range("D34").Value = "5.84"
$ws.Range("E34").Value = "  -3.88%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.128"
$ws.Range("E35").Value = "  -4.33%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  +10.50%  "
$ws.Range("D37").Value = "4.46"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("D39").Value = "3.65"
$ws.Range("E39").Value = "  -5.32%  "
$ws.Range("D40").Value = "2.65"
$ws.Range("E40").Value = "  -8.33%  "
$ws.Range("D41").Value = "98.27"
$ws.Range("E41").Value = "  +7.49%  "
$ws.Range("D42").Value = "1.45"
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("D43").Value = "69.47"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "0.228"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "11.97"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").Value = "110.61"
$ws.Range("E47").Value = "  -6.40%  "
$ws.Range("D48").Value = "5.34"
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("D49").Value = "8.94"
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").Value = "72.53"
$ws.Range("E50").Value = "  +4.38%  "
$ws.Range("E51").Value = "  -1.37%  "
